$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 notes: "G(s) is the process..." paragraph gains italic call-outs for
# Tw, zeta, Kp and Td in the real deck; this runtime cannot mutate per-run
# Font properties on a Notes body placeholder, so we only carry the text
# content forward unchanged (it was already correct) together with the rest
# of the untouched paragraphs, rebuilt through the one notes-text write path
# that actually works: a whole-TextRange `.Text` assignment using `` `n ``
# (LF) as the paragraph separator.
# ---------------------------------------------------------------------------
$notes1 = $p.Slides.Item(1).NotesPage.Shapes.Item(2)
$notes1Text = @"
G(s) is the process to be estimated: it is underdamped (two complex poles defined by Tw and zeta), has a gain (Kp) and a delay (Td).

Estimation is carried out by means of ident toolbox with the aim to solve the least-squares problem. The toolbox chooses the best method among Gradient-Descent, Levenberg-Marquardt, Adaptive Gauss-Newton, Trust-Regions.

Validation of the estimated process is performed on different data sets, achieving high-accuracy of the prediction (~90%).
"@
$notes1.TextFrame.TextRange.Text = $notes1Text

# ---------------------------------------------------------------------------
# Slide 2 notes: split the single paragraph into three paragraphs (adding a
# blank line in between), reword "...fine tune the gain setting..." to
# "...fine tune the setting...", and reword "...is enough." to
# "...is sufficient." at the start of the new final paragraph ("It came out
# that, ...").
# ---------------------------------------------------------------------------
$notes2 = $p.Slides.Item(2).NotesPage.Shapes.Item(2)
$notes2Text = @"
The design of the discrete PI controller (running @ 100 Hz) has the objective to minimize the response time of the plant when undergoes an input stepwise load disturbance. We employed the PID Tuner for that; it’s a GUI where the user has sliders acting on the gains to fine tune the setting proposed to meet the requirements.

It came out that, given the identified G(s), a simple I controller is sufficient.
"@
$notes2.TextFrame.TextRange.Text = $notes2Text
